$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the middle data row (row 17: "YONNATAN ALBERTO GONZALEZ OLIVERA"),
# which shifts the last data row (formerly row 18) up to row 17, and shifts
# the trailing signature rows (23,24) up to (22,23).
$ws.Rows("17").Delete()

# Update the "Valor Mora" total for the remaining bottom data row (now row 17)
$ws.Range("G17").Value = 1423500

# Update the summary totals
$ws.Range("E11").Value = 75108
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2

# Adjust column D width to match new best-fit width after the name column
# content changed (shorter name remaining visible)
$ws.Columns("D").ColumnWidth = 32.584
